$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "No Blink Detected"

# Row 5
$ws.Range("C5").Value = 2

# Row 6
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "Match"

# Row 7
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = "No Blink Detected"

# Row 9
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "No Match"

# Row 10
$ws.Range("C10").Value = 2

# Row 11
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = "Match"
